# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F5"  = 839
    "F6"  = 10
    "F7"  = 295
    "F8"  = 7831
    "F12" = 104
    "F13" = 5
    "F19" = 683
    "F20" = 20
    "F21" = 75
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
